$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2863668.32
$ws.Range("C7").Value = -35.54764866479548
$ws.Range("D7").Value = 2881
$ws.Range("E7").Value = 2881
$ws.Range("F7").Value = 993.9841443943075
$ws.Range("G7").Value = 5.951522361516348
